# Edit: replace the placeholder "Fea/fea/feo/Feo" entries with the tail of
# the insult list shifting up, and append the new "pendejo..perra" words,
# matching the net effect of removing 4 shared strings and adding 9 new
# ones to xl/sharedStrings.xml while leaving rows 1-180 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final text for rows 181-212 (row 180 "Zullenco" and everything above is
# unchanged).
$newWords = @(
    "Tonto", "Tonta", "Bobo", "Boba", "bobo", "boba", "baboso", "babosa",
    "estupido", "estipida", "Estupida", "Estupido", "Puto", "Puta", "Putazo",
    "Chingado", "Chingada", "Idiota", "idiota", "mamador", "mamadora",
    "imbecil", "Imbecil", "pendejo", "pendeja", "Pendejazo", "Pendejo",
    "puta", "puto", "putisimo", "putisima", "perra"
)

$startRow = 181
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value2 = $newWords[$i]
}

# Rows 185 through the new end (212) pick up the same style already used by
# A186 ("s=1" in the original workbook) instead of the plain A2-style
# formatting ("s=2") the rest of the list uses.
$fmtSource = $ws.Range("A186")
$fmtSource.Copy()
$fmtTarget = $ws.Range("A185:A212")
$fmtTarget.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
